# Refresh the cryptos price/volume table with the latest scrape
# (scheduled GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.616.07'
$ws.Range('E2').Value = '  -2.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.665.72'
$ws.Range('E3').Value = '  -4.02%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.65'
$ws.Range('E5').Value = '  -2.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.510'
$ws.Range('E6').Value = '  -2.64%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.06'
$ws.Range('E8').Value = '  -0.39%  '
$ws.Range('E9').Value = '  -1.31%  '
$ws.Range('E10').Value = '  -2.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0879'
$ws.Range('E11').Value = '  -2.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.901.26'
$ws.Range('E12').Value = '  -4.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.646.85'
$ws.Range('E13').Value = '  -5.27%  '
$ws.Range('E14').Value = '  -3.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.570'
$ws.Range('E15').Value = '  +1.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.37'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.599.71'
$ws.Range('E17').Value = '  -2.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '241.79'
$ws.Range('E18').Value = '  -0.56%  '
$ws.Range('E19').Value = '  -3.65%  '
$ws.Range('E20').Value = '  -4.51%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('E22').Value = '  -3.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.36'
$ws.Range('E23').Value = '  -3.77%  '
$ws.Range('E24').Value = '  -3.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.74'
$ws.Range('E25').Value = '  -2.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.27'
$ws.Range('E26').Value = '  -3.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.37'
$ws.Range('E27').Value = '  -2.12%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('E29').Value = '  -2.66%  '
$ws.Range('E30').Value = '  +2.13%  '
$ws.Range('E31').Value = '  -2.26%  '
$ws.Range('E32').Value = '  -2.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.461.55'
$ws.Range('E33').Value = '  -3.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.13'
$ws.Range('E34').Value = '  -4.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.59'
$ws.Range('E35').Value = '  -4.32%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.931'
$ws.Range('E36').Value = '  -3.76%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.37'
$ws.Range('E37').Value = '  -1.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.578'
$ws.Range('E38').Value = '  -5.00%  '
$ws.Range('E39').Value = '  -2.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '69.91'
$ws.Range('E40').Value = '  -1.21%  '
$ws.Range('E41').Value = '  -4.00%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.22'
$ws.Range('E43').Value = '  -3.78%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.41'
$ws.Range('E44').Value = '  -5.44%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.793'
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.809.17'
$ws.Range('E46').Value = '  -4.02%  '
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.83'
$ws.Range('E48').Value = '  -2.59%  '
$ws.Range('E49').Value = '  -6.62%  '
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.91'
$ws.Range('E51').Value = '  -3.97%  '
